$d = $word.ActiveDocument

# Remove the "ps sign.jpg" signature picture first (it sits in the
# paragraph right after the employee-details table). Deleting the
# InlineShape before the table keeps story/paragraph indices valid.
$sig = $d.InlineShapes.Item(1)
[void]$sig.Delete()

# Remove the employee DAR/NOC details table (the second table in the
# document - the first one is the header/letterhead table).
$tbl = $d.Tables.Item(2)
[void]$tbl.Delete()
